$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.332878589630127
$ws.Range("B1").Value = 2.132723569869995
$ws.Range("C1").Value = 4.827506542205811
$ws.Range("D1").Value = 3.432815313339233
$ws.Range("E1").Value = 1.309719562530518
